{"js": "// Update the \"Features of the transfer\" table:\n//   1. Rename the \"Who initiated this academy transfer\" row label/placeholder\n//      to \"Reason for this transfer\" / \"ReasonForTransfer\".\n//   2. Remove the now-redundant \"Has this transfer started because of an\n//      intervention with the academy or trust?\" row and the\n//      \"More details about the transfer\" row that followed it.\n\nconst body = context.document.body;\n\n// 1. Update the visible label text for the first row of that table.\nconst labelResults = body.search(\"Who initiated this academy transfer\", { matchCase: true });\nlabelResults.load(\"items\");\nawait context.sync();\nlabelResults.items[0].insertText(\"Reason for this transfer\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. Update the bracketed placeholder token in the same row's answer cell.\nconst placeholderResults = body.search(\"WhoInitiatedTheTransfer\", { matchCase: true });\nplaceholderResults.load(\"items\");\nawait context.sync();\nplaceholderResults.items[0].insertText(\"ReasonForTransfer\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. Locate the table that contains this row (\"Features of the transfer\")\n//    and delete the two rows that followed it (now obsolete).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet targetTable = null;\nfor (let i = 0; i < tables.items.length; i++) {\n  const t = tables.items[i];\n  t.load(\"values\");\n}\nawait context.sync();\n\nfor (let i = 0; i < tables.items.length; i++) {\n  const t = tables.items[i];\n  if (t.values && t.values.length && t.values[0][0] === \"Reason for this transfer\") {\n    targetTable = t;\n    break;\n  }\n}\n\nif (targetTable) {\n  targetTable.rows.load(\"items\");\n  await context.sync();\n  // Delete from the bottom up so earlier indices stay valid.\n  // Row 1 = \"Has this transfer started...\" ; Row 2 = \"More details about the transfer\"\n  targetTable.rows.items[2].delete();\n  targetTable.rows.items[1].delete();\n  await context.sync();\n}\n", "ps1": "# Update the \"Features of the transfer\" table:\n#   1. Rename the \"Who initiated this academy transfer\" row label/placeholder\n#      to \"Reason for this transfer\" / \"ReasonForTransfer\".\n#   2. Remove the now-redundant \"Has this transfer started because of an\n#      intervention with the academy or trust?\" row and the\n#      \"More details about the transfer\" row that followed it.\n\n$d = $word.ActiveDocument\n\n# 1. Update the visible label text for the first row of that table.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Who initiated this academy transfer\"\n$find.Replacement.Text = \"Reason for this transfer\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2. Update the bracketed placeholder token in the same row's answer cell.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"WhoInitiatedTheTransfer\"\n$find2.Replacement.Text = \"ReasonForTransfer\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n# 3. Locate the table that contains this row (\"Features of the transfer\")\n#    and delete the two rows that followed it (now obsolete).\nforeach ($t in $d.Tables) {\n    if ($t.Rows.Count -ge 3 -and $t.Cell(1,1).Range.Text.Contains(\"Reason for this transfer\")) {\n        $t.Rows.Item(3).Delete()\n        $t.Rows.Item(2).Delete()\n        break\n    }\n}\n"}
